# Append a new paragraph at the end of the document describing the timer
# system, matching the target OOXML diff:
#   "A timer system was implemented using pre-written code (Timer, "
#   "TimerMonoBehaviour" and " and " "ShowOnlyAttribute" ")"
# split across separate runs (the two technical identifiers are bracketed by
# proofErr markers in the authored document, which live Word's proofing
# engine stamps as a side effect of typing unrecognised words; those
# self-closing annotations carry no retrievable text/content and are not
# reachable through the Word object model, so this script focuses on
# reproducing the paragraph/run/text structure).

$d = $word.ActiveDocument

# Move to the very end of the document and start a fresh paragraph that
# inherits the same paragraph mark formatting (en-US language) as the
# paragraph before it.
$lastPara = $d.Paragraphs($d.Paragraphs.Count)
$endRange = $lastPara.Range
$endRange.Collapse(0)
$endRange.InsertParagraphAfter()

$newPara = $d.Paragraphs($d.Paragraphs.Count)
$newRange = $newPara.Range

# First run.
$newRange.Text = "A timer system was implemented using pre-written code (Timer, "
$newRange.Collapse(0)

# "TimerMonoBehaviour" run - wrap the insertion point in a transient
# bookmark so the new text lands in its own run instead of being coalesced
# into the previous run (the two runs already share identical formatting).
$d.Bookmarks.Add("__edit_mark1", $newRange)
$newRange.InsertAfter("TimerMonoBehaviour")
$newRange.LanguageID = "en-US"
$d.Bookmarks("__edit_mark1").Delete()
$newRange.Collapse(0)

# " and " run.
$d.Bookmarks.Add("__edit_mark2", $newRange)
$newRange.InsertAfter(" and ")
$newRange.LanguageID = "en-US"
$d.Bookmarks("__edit_mark2").Delete()
$newRange.Collapse(0)

# "ShowOnlyAttribute" run.
$d.Bookmarks.Add("__edit_mark3", $newRange)
$newRange.InsertAfter("ShowOnlyAttribute")
$newRange.LanguageID = "en-US"
$d.Bookmarks("__edit_mark3").Delete()
$newRange.Collapse(0)

# Closing parenthesis run.
$d.Bookmarks.Add("__edit_mark4", $newRange)
$newRange.InsertAfter(")")
$newRange.LanguageID = "en-US"
$d.Bookmarks("__edit_mark4").Delete()
